$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.113.05"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "1.823.72"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.37%  "

$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07305"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8708"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "1.871.22"
$ws.Range("E12").Value = "  +0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07594"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.352"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.490"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008640"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.44%  "

$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "27.396.29"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("E21").Value = "  -2.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.205"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.60%  "

$ws.Range("E23").Value = "  -1.73%  "

$ws.Range("D24").Value = "2.092.24"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.874"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.096"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.47%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.093"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.06%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08917"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.951"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7347"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.93%  "

$ws.Range("E34").Value = "  -2.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.66%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.483"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.073"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05255"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("E40").Value = "  -2.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.927"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.143"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5205"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1629"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.279"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4885"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.009"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06254"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
